$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D (MAE) before the current "Tipo" column,
# shifting the existing "Tipo"/"single" data into column E.
$ws.Columns.Item(4).Insert()

# New header for the inserted column D
$ws.Range("D1").Value = "MAE"

# Copy the header formatting (bold font, border, centered alignment) from an
# existing header cell onto the new header cell.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# New MAE value for row 2
$ws.Range("D2").Value = 0.4938974175076344

# Update existing MSE (B2) and R2 (C2) values
$ws.Range("B2").Value = 0.4140520323401509
$ws.Range("C2").Value = 0.9918354608100279
